# Fruta / hortaliza, semanal
# Insert a new weekly record at row 256 (pushing existing rows 256-279 down
# to 257-280) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 256; this shifts rows 256:279
# down to 257:280 and carries the date-format style from row 256 (s="2")
# onto the new D256, matching the surrounding rows.
$ws.Rows(256).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A256").Value = 1
$ws.Range("B256").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C256").Value = "Arica y Parinacota"
$ws.Range("D256").Value = 45194
$ws.Range("E256").Value = 15
$ws.Range("F256").Value = 100114001
$ws.Range("G256").Value = "Papa"
$ws.Range("H256").Value = "Asterix"
$ws.Range("I256").Value = "1a (cosecha)"
$ws.Range("J256").Value = 1000
$ws.Range("K256").Value = 30000
$ws.Range("L256").Value = 31000
$ws.Range("M256").Value = 30500
$ws.Range("N256").Value = "$/saco 25 kilos"
$ws.Range("O256").Value = "Provincia de Melipilla"
$ws.Range("P256").Value = 1220
$ws.Range("Q256").Value = 25
$ws.Range("R256").Value = "Hortaliza"
